$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 262, shifting rows 262-316 down to 263-317
$ws.Rows.Item(262).Insert()

# Populate the newly inserted row 262 with the new data
$ws.Cells.Item(262, 1).Value = 3
$ws.Cells.Item(262, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(262, 3).Value = "Coquimbo"
$ws.Cells.Item(262, 4).Value = 44637
$ws.Cells.Item(262, 5).Value = 5
$ws.Cells.Item(262, 6).Value = 100112043
$ws.Cells.Item(262, 7).Value = "Pepino ensalada"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 105
$ws.Cells.Item(262, 11).Value = 18000
$ws.Cells.Item(262, 12).Value = 18500
$ws.Cells.Item(262, 13).Value = 18262
$ws.Cells.Item(262, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(262, 15).Value = "Limache"
$ws.Cells.Item(262, 16).Value = 261
$ws.Cells.Item(262, 17).Value = 70
$ws.Cells.Item(262, 18).Value = "Hortaliza"
